$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "params": update the event table, drop the trailing row, extend the
# comment column with two new entries. Cells are touched row-by-row (date
# column first, then the numeric columns, then the comment column) so that
# newly introduced shared strings land in the same order a person editing
# top-to-bottom would create them.
# ---------------------------------------------------------------------------
$params = $wb.Worksheets.Item("params")

# Delete the last data row (15.09.2020 / "Fin supposée ..." row) - the table
# now ends at row 13 instead of row 14.
$params.Rows.Item(14).Delete()

# Row 2 - 18.01.2020 (unchanged dates/comment, values already correct)
$params.Range("A2").Value = "18.01.2020"
$params.Range("B2").Value = 1.2
$params.Range("C2").Value = 0.2
$params.Range("D2").Value = 0.2
$params.Range("L2").Value = "Premier patient dans fichier de données"

# Row 3 - 25.02.2020
$params.Range("A3").Value = "25.02.2020"
$params.Range("B3").Value = 1.2
$params.Range("C3").Value = 0.2
$params.Range("D3").Value = 0.2
$params.Range("L3").Value = "Début épidémie dans le canton de VD"

# Row 4 - 28.02.2020 (new date, new comment)
$params.Range("A4").Value = "28.02.2020"
$params.Range("B4").Value = 1.19
$params.Range("C4").Value = 0.2
$params.Range("D4").Value = 0.2
$params.Range("L4").Value = "Interdiction rassemblements > 1000 personnes"

# Row 5 - 13.03.2020 (new date, new comment)
$params.Range("A5").Value = "13.03.2020"
$params.Range("B5").Value = 1.18
$params.Range("C5").Value = 0.2
$params.Range("D5").Value = 0.2
$params.Range("L5").Value = "Interdiction rassemblements > 100 personnes, fermeture écoles, contrôles aux frontières"

# Row 6 - 16.03.2020
$params.Range("A6").Value = "16.03.2020"
$params.Range("B6").Value = 1.17
$params.Range("C6").Value = 0.2
$params.Range("D6").Value = 0.2
$params.Range("L6").Value = "Situation extraordinaire : fermeture des commerces non essentiels, fermeture partielle des frontières"

# Row 7 - 28.03.2020
$params.Range("A7").Value = "28.03.2020"
$params.Range("B7").Value = 1.05
$params.Range("C7").Value = 0.2
$params.Range("D7").Value = 0.18
$params.Range("L7").Value = "Début de l’effet du confinement (megp → 1)"

# Row 8 - 05.04.2020
$params.Range("A8").Value = "05.04.2020"
$params.Range("B8").Value = 1.03
$params.Range("C8").Value = 0.2
$params.Range("D8").Value = 0.18

# Row 9 - 10.04.2020
$params.Range("A9").Value = "10.04.2020"
$params.Range("B9").Value = 1.02
$params.Range("C9").Value = 0.2
$params.Range("D9").Value = 0.19

# Row 10 - 15.04.2020 (new date)
$params.Range("A10").Value = "15.04.2020"
$params.Range("B10").Value = 1.01
$params.Range("C10").Value = 0.2
$params.Range("D10").Value = 0.19
$params.Range("L10").ClearContents()

# Row 11 - 01.05.2020 (new date)
$params.Range("A11").Value = "01.05.2020"
$params.Range("B11").Value = 1.005
$params.Range("C11").Value = 0.2
$params.Range("D11").Value = 0.2

# Row 12 - 15.05.2020 (new date)
$params.Range("A12").Value = "15.05.2020"
$params.Range("B12").Value = 1.002
$params.Range("C12").Value = 0.2
$params.Range("D12").Value = 0.2
$params.Range("L12").ClearContents()

# Row 13 - 01.06.2020 (new date, new comment: now the "end of epidemic" row)
$params.Range("A13").Value = "01.06.2020"
$params.Range("B13").Value = 1
$params.Range("C13").Value = 0.2
$params.Range("D13").Value = 0.2
$params.Range("L13").Value = "Fin supposée de l’épidémie (megp=1)"

# ---------------------------------------------------------------------------
# Sheet "age_distrib": new age buckets + refreshed survival shares.
# ---------------------------------------------------------------------------
$age = $wb.Worksheets.Item("age_distrib")
$age.Range("B1").Value = "0-69"
$age.Range("C1").Value = "70-84"
$age.Range("D1").Value = "85-119"
$age.Range("B2").Value = 0.494
$age.Range("C2").Value = 0.332
$age.Range("D2").Value = 0.174

# ---------------------------------------------------------------------------
# Sheet "sex_distrib": same new age buckets + refreshed survival shares.
# ---------------------------------------------------------------------------
$sex = $wb.Worksheets.Item("sex_distrib")
$sex.Range("B1").Value = "0-69"
$sex.Range("C1").Value = "70-84"
$sex.Range("D1").Value = "85-119"
$sex.Range("B2").Value = 0.372
$sex.Range("C2").Value = 0.411
$sex.Range("D2").Value = 0.528
